function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.079.32"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.442.88"
$ws.Range("E3").Value = "  +0.24%  "
Set-TextValue $ws "D4" "0.998"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue $ws "D5" "579.93"
$ws.Range("E5").Value = "  +2.05%  "
Set-TextValue $ws "D6" "143.16"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  -0.01%  "
Set-TextValue $ws "D8" "0.530"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "2.437.64"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("E12").Value = "  -0.95%  "
Set-TextValue $ws "D13" "0.344"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("E14").Value = "  -1.72%  "
Set-TextValue $ws "D15" "0.0000172"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").Value = "2.800.52"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").Value = "62.069.12"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "2.425.16"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  -3.60%  "
Set-TextValue $ws "D20" "7.15"
$ws.Range("E20").Value = "  -1.63%  "
Set-TextValue $ws "D21" "327.83"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("E24").Value = "  +0.15%  "
Set-TextValue $ws "D25" "65.57"
Set-TextValue $ws "D26" "9.34"
$ws.Range("E26").Value = "  +6.73%  "
Set-TextValue $ws "D27" "606.95"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "2.565.83"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "0.0₃0945"
$ws.Range("E30").Value = "  -6.27%  "
$ws.Range("E31").Value = "  -4.85%  "
Set-TextValue $ws "D32" "7.97"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D33" "1.88"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D34" "0.140"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("E35").Value = "  -4.40%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws "D37" "0.375"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D38" "1.42"
$ws.Range("E38").Value = "  -5.58%  "
Set-TextValue $ws "D39" "149.08"
$ws.Range("E39").Value = "  +2.81%  "
Set-TextValue $ws "D40" "5.29"
$ws.Range("E40").Value = "  -0.33%  "
Set-TextValue $ws "D41" "18.32"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("E42").Value = "  -2.76%  "
Set-TextValue $ws "D43" "42.50"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("E44").Value = "  -0.01%  "
Set-TextValue $ws "D45" "2.44"
$ws.Range("E45").Value = "  -5.33%  "
Set-TextValue $ws "D46" "142.41"
$ws.Range("E46").Value = "  -3.46%  "
Set-TextValue $ws "D47" "3.62"
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("E48").Value = "  +1.09%  "
Set-TextValue $ws "D49" "0.0522"
$ws.Range("E49").Value = "  -1.23%  "
Set-TextValue $ws "D50" "19.40"
$ws.Range("E50").Value = "  -6.75%  "
$ws.Range("D51").Value = "0.0₆0233"
$ws.Range("E51").Value = "  +8.44%  "

